$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.843.39"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "2.438.96"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'559.43"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").Value = "'162.42"
$ws.Range("E6").Value = "  -1.30%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'0.510"
$ws.Range("E8").Value = "  -0.54%  "

$ws.Range("E9").Value = "  +9.75%  "

$ws.Range("E10").Value = "  -2.16%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("E12").Value = "  -4.86%  "

$ws.Range("E13").Value = "  +4.86%  "

$ws.Range("D14").Value = "68.715.95"
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("D15").Value = "2.886.52"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").Value = "'23.30"
$ws.Range("E16").Value = "  -1.35%  "

$ws.Range("D17").Value = "2.440.12"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").Value = "'339.44"
$ws.Range("E19").Value = "  +0.22%  "

$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("E21").Value = "  +0.92%  "

$ws.Range("E22").Value = "  +3.67%  "

$ws.Range("D24").Value = "'66.33"
$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("E25").Value = "  +1.61%  "

$ws.Range("D26").Value = "2.567.54"
$ws.Range("E26").Value = "  -1.13%  "

$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("D28").Value = "'0.990"
$ws.Range("E28").Value = "  -0.88%  "

$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +2.86%  "

$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").Value = "'429.62"
$ws.Range("E33").Value = "  -0.18%  "

$ws.Range("E34").Value = "  -1.79%  "

$ws.Range("D35").Value = "'160.19"
$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("E38").Value = "  +0.82%  "

$ws.Range("E39").Value = "  -2.35%  "

$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("E41").Value = "  +2.24%  "

$ws.Range("D42").Value = "'4.37"
$ws.Range("E42").Value = "  -1.83%  "

$ws.Range("D43").Value = "'1.07"
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("D44").Value = "'2.07"
$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("D46").Value = "'130.33"
$ws.Range("E46").Value = "  -0.39%  "

$ws.Range("D47").Value = "'0.0720"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("D48").Value = "'0.483"
$ws.Range("E48").Value = "  -1.00%  "

$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("E50").Value = "  +2.92%  "

$ws.Range("D51").Value = "'0.0921"
$ws.Range("E51").Value = "  +0.07%  "
